$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right after the header (row 1), pushing the
# existing data down by two rows.
$ws.Rows.Item(2).Resize(2).Insert()

# Populate the two newly inserted rows with the new sample readings.
$ws.Range("A2").Value = -3.823432922363281
$ws.Range("B2").Value = 5.642979621887207
$ws.Range("C2").Value = 4.475735664367676

$ws.Range("A3").Value = -4.128349304199219
$ws.Range("B3").Value = 5.92755126953125
$ws.Range("C3").Value = 5.054780960083008

# The trailing three rows of the original data (now shifted to rows
# 22-24) are dropped so the sheet ends at row 21.
$ws.Rows.Item("22:24").Delete()
